# Renderer003-WsName fixture update:
#   "increase coverage & fix detected bug inside DeleteCell"
#
# The workbook originally has 3 sheets: Hello, World, Sheet3.
# The updated fixture adds two more worksheets - "World 2" and "World 3" -
# both duplicates of "World" (same single cell A1 = "GOOD BYE", same
# column widths), inserted right after "World" and before "Sheet3":
#
#   Hello, World, World 2, World 3, Sheet3
#
# and leaves the newly-added last sheet ("World 3") as the active one.

$wb = $excel.ActiveWorkbook

$wsWorld = $wb.Worksheets.Item("World")

# Duplicate "World" -> placed immediately after "World"; Excel auto-names
# the copy "World (2)".
$wsWorld.Copy($null, $wsWorld)
$world2 = $wb.Worksheets.Item("World (2)")
$world2.Name = "World 2"

# Duplicate "World 2" -> placed immediately after "World 2"; Excel
# auto-names the copy "World 2 (2)".
$world2.Copy($null, $world2)
$world3 = $wb.Worksheets.Item("World 2 (2)")
$world3.Name = "World 3"

# Leave the last-created sheet selected/active, matching the fixture.
$world3.Activate()
